$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at the top of the data block (becomes row 38),
# shifting the existing rows 38:79 down to 39:80.
$ws.Rows("38:38").Insert()
$ws.Cells.Item(38,1).Value = 5
$ws.Cells.Item(38,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(38,3).Value = "Maule"
$ws.Cells.Item(38,4).Value = "2021-12-24"
$ws.Cells.Item(38,5).Value = 7
$ws.Cells.Item(38,6).Value = 100112030
$ws.Cells.Item(38,7).Value = "Poroto granado"
$ws.Cells.Item(38,8).Value = "Sin especificar"
$ws.Cells.Item(38,9).Value = "Primera"
$ws.Cells.Item(38,10).Value = 200
$ws.Cells.Item(38,11).Value = 35000
$ws.Cells.Item(38,12).Value = 35000
$ws.Cells.Item(38,13).Value = 35000
$ws.Cells.Item(38,14).Value = "$/saco 25 kilos"
$ws.Cells.Item(38,15).Value = "Región del Maule"
$ws.Cells.Item(38,16).Value = 1400
$ws.Cells.Item(38,17).Value = 25
$ws.Cells.Item(38,18).Value = "Hortaliza"

# Insert a second new data row further down (becomes row 57 after the
# previous insert shifted everything), shifting rows 57:80 down to 58:81.
$ws.Rows("57:57").Insert()
$ws.Cells.Item(57,1).Value = 5
$ws.Cells.Item(57,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(57,3).Value = "Maule"
$ws.Cells.Item(57,4).Value = "2021-12-23"
$ws.Cells.Item(57,5).Value = 7
$ws.Cells.Item(57,6).Value = 100112030
$ws.Cells.Item(57,7).Value = "Poroto granado"
$ws.Cells.Item(57,8).Value = "Sin especificar"
$ws.Cells.Item(57,9).Value = "Primera"
$ws.Cells.Item(57,10).Value = 200
$ws.Cells.Item(57,11).Value = 35000
$ws.Cells.Item(57,12).Value = 35000
$ws.Cells.Item(57,13).Value = 35000
$ws.Cells.Item(57,14).Value = "$/saco 25 kilos"
$ws.Cells.Item(57,15).Value = "Región del Maule"
$ws.Cells.Item(57,16).Value = 1400
$ws.Cells.Item(57,17).Value = 25
$ws.Cells.Item(57,18).Value = "Hortaliza"
